$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction6")

$values = @(1,2,3,5,6,7,8,10,12,13,14,15,17,18,19,20,21,22,24,27,28,29,30,31,32,33,34,36,38,39,40,41,42,44,46,47,49,51,52,53,54,55,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,73,74)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
